$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Row 2: Remarks text is unchanged content-wise (shared string just gets
#     renumbered internally when the table grows) - no explicit edit needed.

# --- Row 5: "Write function to create schema for all devices" checklist -> done,
#     remarks text appended with ", Done for MDM"
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = "Function will hold list of dictionaries of devices, Done for MDM"
$ws.Range("A5:D5").EntireRow.RowHeight = 28.8

# --- Row 6: checklist -> done, remarks now "Done"
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = "Done"

# --- Row 7: rewritten post-API call description, checklist -> done
$ws.Range("B7").Value = "Write the call to run the post API"
$ws.Range("B7").WrapText = $true
$ws.Range("C7").Value = $true

# --- Rows 8-9 shift up (old row 7/8 content moves down one slot since the
#     "join lists" task got replaced by the new row 7 call-writing task)
$ws.Range("B8").Value = "Write snippet to join 2 lists in serial order"
$ws.Range("B9").Value = "Test code to check if it works"

# --- New row 10: add the task that used to live at row 9. Insert a
#     formatted row first (inherits formatting from row 9 above it), then
#     register it with the table so Table1's range expands from A1:D9 to
#     A1:D10, then fill in the values.
$ws.Rows.Item(10).Insert()
$newRow = $lo.ListRows.Add()
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Modify the existing codes to run post API"
$ws.Range("C10").Value = $false

# --- Selection moves to C7 (matches the saved cursor position in the file)
$ws.Range("C7").Select()
